$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Startseite")

# Delete the 5 empty placeholder rows (rows 2-6), shifting the data rows
# below them up so the table becomes contiguous starting at row 2.
$ws.Rows("2:6").Delete()

# Move selection to reflect where the user ended up after the edit.
$ws.Range("B25").Select()
